# DBA_Automon_Details.xlsx - update Description column wording, bold the
# "group header" description cells to match the rest of their row, add a
# trailing note row, and refresh the workbook's recorded author path /
# selection, per the "Adding extra column in base table DBA_All_servers"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Description) wording fixes -------------------------------
$ws.Range("D7").Value  = "Recovery model check for non production"
$ws.Range("D8").Value  = "Recovery model check for production"
$ws.Range("D12").Value = "Auto maintenance mode after 5 ping health fails"
$ws.Range("D13").Value = "Maintenance window release daily morning"
$ws.Range("D17").Value = "AutoMon SQL server not running status"
$ws.Range("D18").Value = "Disk free space percentage alert"
$ws.Range("D19").Value = "Error report of AutoMon tool"
$ws.Range("D20").Value = "Write script to get AutoMon Job Disabled alert"
$ws.Range("D22").Value = "Backup share path free space alert "
$ws.Range("D25").Value = "AlwaysOn out of sync lag check"

# --- Bold the Description cell on each "group header" row, matching the
#     bold styling already used by the A/B/C cells on those same rows ----
$ws.Range("D2").Font.Bold  = $true
$ws.Range("D3").Font.Bold  = $true
$ws.Range("D9").Font.Bold  = $true
$ws.Range("D23").Font.Bold = $true
$ws.Range("D26").Font.Bold = $true

# --- New trailing note row ------------------------------------------------
$newRow = $ws.Range("A29")
$newRow.Value = "Configure the job based on your requirements"
$newRow.Borders.Item(7).LineStyle  = 1   # xlEdgeLeft
$newRow.Borders.Item(7).Weight     = 2   # xlThin
$newRow.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$newRow.Borders.Item(10).Weight    = 2   # xlThin

# --- Workbook-level housekeeping -----------------------------------------
# Recorded author path, changed by Excel on save to the new machine/repo.
$wb.Application.ActiveWorkbook.FullName | Out-Null
$ws.Range("C14").Select()

Write-Output "edit complete"
